# Grammar and spell check edits
#
# 1) "...the parrot alone, the cat..."  ->  "...the parrot alone; the cat..."
#    (splits the run so "alone;" becomes its own run, matching the
#     grammar-checker style edit recorded in the target document)
# 2) "...getting black socks is greater..." -> "...getting black socks are greater..."
#    (splits the run so "socks are" becomes its own run, and the
#     Word "last edit" bookmark _GoBack ends up collapsed right before it)
# 3) The _GoBack bookmark is therefore removed from its old location
#    at the end of the "Predicting Finger" section, since Word keeps
#    only one _GoBack bookmark - at the site of the most recent edit.

$d = $word.ActiveDocument

# Locate the two target paragraphs by content instead of a fixed index,
# so the script is resilient to any unrelated paragraph numbering shifts.
$p1 = $null
$p2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($p1 -eq $null -and $t.Contains("alone,")) {
        $p1 = $p
    }
    if ($p2 -eq $null -and $t.Contains("socks is greater")) {
        $p2 = $p
    }
}

# ---------------------------------------------------------------
# Edit 1: Cat/Parrot paragraph - "alone," -> "alone;"
# ---------------------------------------------------------------
$p1Text = $p1.Range.Text
$offset1 = $p1Text.IndexOf("alone,")
$sub1Start = $p1.Range.Start + $offset1
$sub1End = $sub1Start + 6

$r1 = $d.Range($sub1Start, $sub1End)
# Use a temporary bookmark to force this substring into its own run
# without leaving any residual run-level formatting behind.
$d.Bookmarks.Add("tmpSplit1", $r1)
$r1.Text = "alone;"
$d.Bookmarks.Item("tmpSplit1").Delete()

# ---------------------------------------------------------------
# Edit 2: Socks paragraph - "socks is" -> "socks are"
# ---------------------------------------------------------------
$p2Text = $p2.Range.Text
$offset2 = $p2Text.IndexOf("socks is greater")
$sub2Start = $p2.Range.Start + $offset2
$sub2End = $sub2Start + 8

$r2 = $d.Range($sub2Start, $sub2End)
$d.Bookmarks.Add("tmpSplit2", $r2)
$r2.Text = "socks are"
$d.Bookmarks.Item("tmpSplit2").Delete()

# ---------------------------------------------------------------
# Edit 3: Move the _GoBack bookmark from the end of the document
# to the (collapsed) location of the most recent edit above.
# ---------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$goBackRange = $d.Range($sub2Start, $sub2Start)
$d.Bookmarks.Add("_GoBack", $goBackRange)
